# NIT-9011266848.xlsx — "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta"
#
# Updates the "Estado de Cuenta" (account statement) for company
# INSTITUTO DE FERTILIDAD HUMANA CARIBE S.A.S. (renamed from ...CARTAGENA):
#   - corrects the company/razon social name
#   - updates the VALOR MORA total
#   - updates the worker count (Cant. Trabajadores) and periods count (Cant. Periodos)
#   - drops the second worker (YURI EIDY VILLADA HOYOS) entirely, keeping only
#     ISOLINA MARIA MARVAL PEREZ, whose 3 remaining period rows are reordered
#     chronologically (2405, 2406, 2407) with the corrected "Valor Mora" amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Razon social: CARTAGENA -> CARIBE
$ws.Range("E7").Value = "INSTITUTO DE FERTILIDAD HUMANA CARIBE S.A.S."

# VALOR MORA (total overdue amount)
$ws.Range("E11").Value = 781387

# Cant. Trabajadores / Cant. Periodos
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 3

# Remove one of ISOLINA's period rows (period 2407, which will be re-added
# below with the last row's special border styling) and both of YURI's rows,
# leaving exactly 3 data rows (16:18) for ISOLINA.
$ws.Rows("16:16").Delete()
$ws.Rows("18:18").Delete()

# Rows 16:18 now hold ISOLINA's 2406/2405 rows plus YURI's last (styled) row.
# Rewrite them in chronological period order with the correct mora values.
$ws.Range("E16").Value = "2405"
$ws.Range("F16").Value = 360640
$ws.Range("G16").Value = 9016000

$ws.Range("E17").Value = "2406"
$ws.Range("F17").Value = 360640
$ws.Range("G17").Value = 9016000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "2000021510"
$ws.Range("D18").Value = "ISOLINA MARIA MARVAL PEREZ"
$ws.Range("E18").Value = "2407"
$ws.Range("F18").Value = 60107
$ws.Range("G18").Value = 9016000
